$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.131.57"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.801.71"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.41%  "

Set-TextValue $ws.Range("D5") "311.15"
$ws.Range("E5").Value = "  -1.27%  "

Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  -2.44%  "

Set-TextValue $ws.Range("D8") "0.3883"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("E9").Value = "  -2.58%  "

Set-TextValue $ws.Range("D10") "1.098"
$ws.Range("E10").Value = "  -0.06%  "

Set-TextValue $ws.Range("D11") "40.91"
$ws.Range("E11").Value = "  -2.18%  "

Set-TextValue $ws.Range("D12") "6.320"
$ws.Range("E12").Value = "  -0.11%  "

Set-TextValue $ws.Range("D13") "1.002"
$ws.Range("E13").Value = "  -0.37%  "

Set-TextValue $ws.Range("D14") "20.25"

$ws.Range("D15").Value = "1.798.35"
$ws.Range("E15").Value = "  -0.65%  "

Set-TextValue $ws.Range("D16") "7.267"
$ws.Range("E16").Value = "  -1.08%  "

Set-TextValue $ws.Range("D17") "91.98"
$ws.Range("E17").Value = "  -0.62%  "

Set-TextValue $ws.Range("D18") "0.00001071"
$ws.Range("E18").Value = "  -1.91%  "

Set-TextValue $ws.Range("D19") "0.06572"
$ws.Range("E19").Value = "  -0.40%  "

Set-TextValue $ws.Range("D20") "1.001"
$ws.Range("E20").Value = "  -0.36%  "

Set-TextValue $ws.Range("D21") "17.23"
$ws.Range("E21").Value = "  -1.11%  "

Set-TextValue $ws.Range("D22") "5.965"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").Value = "28.135.23"
$ws.Range("E23").Value = "  -0.61%  "

Set-TextValue $ws.Range("D24") "11.06"
$ws.Range("E24").Value = "  -1.05%  "

Set-TextValue $ws.Range("D25") "2.249"
$ws.Range("E25").Value = "  +0.51%  "

Set-TextValue $ws.Range("D26") "160.93"
$ws.Range("E26").Value = "  +2.14%  "

Set-TextValue $ws.Range("D27") "2.428"
$ws.Range("E27").Value = "  +1.21%  "

$ws.Range("D28").Value = "2.005.65"
$ws.Range("E28").Value = "  -0.44%  "

Set-TextValue $ws.Range("D29") "20.26"
$ws.Range("E29").Value = "  -1.30%  "

Set-TextValue $ws.Range("D30") "126.93"
$ws.Range("E30").Value = "  +3.06%  "

Set-TextValue $ws.Range("D31") "0.1089"
$ws.Range("E31").Value = "  -1.41%  "

Set-TextValue $ws.Range("D32") "1.048"
$ws.Range("E32").Value = "  -1.25%  "

Set-TextValue $ws.Range("D33") "3.646"
$ws.Range("E33").Value = "  -0.64%  "

Set-TextValue $ws.Range("D34") "5.524"
$ws.Range("E34").Value = "  -1.09%  "

Set-TextValue $ws.Range("D35") "0.07026"
$ws.Range("E35").Value = "  -2.45%  "

Set-TextValue $ws.Range("D36") "9.044"
$ws.Range("E36").Value = "  +3.75%  "

Set-TextValue $ws.Range("D37") "0.02341"
$ws.Range("E37").Value = "  +1.12%  "

Set-TextValue $ws.Range("D38") "0.2156"
$ws.Range("E38").Value = "  -0.56%  "

Set-TextValue $ws.Range("D39") "5.013"
$ws.Range("E39").Value = "  -0.61%  "

Set-TextValue $ws.Range("D40") "11.47"
$ws.Range("E40").Value = "  -5.38%  "

Set-TextValue $ws.Range("D41") "0.6113"
$ws.Range("E41").Value = "  -1.33%  "

Set-TextValue $ws.Range("D42") "1.001"
$ws.Range("E42").Value = "  -0.39%  "

Set-TextValue $ws.Range("D43") "1.155"
$ws.Range("E43").Value = "  -1.34%  "

Set-TextValue $ws.Range("D44") "13.10"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.5905"
$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D46") "1.294"
$ws.Range("E46").Value = "  -6.27%  "

$ws.Range("E47").Value = "  -1.19%  "

Set-TextValue $ws.Range("D48") "125.01"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("E49").Value = "  -0.65%  "

Set-TextValue $ws.Range("D50") "1.900"
$ws.Range("E50").Value = "  -1.66%  "

Set-TextValue $ws.Range("D51") "0.06740"
$ws.Range("E51").Value = "  -1.30%  "
